$wb = $excel.ActiveWorkbook

# --- MassInertia: update roll-axis height value (I2) and selection ---
$wsMass = $wb.Worksheets.Item("MassInertia")
$wsMass.Range("I2").Value = 0.43
$wsMass.Activate()
$wsMass.Range("I3").Select()

# --- TireRear: update selected cell only ---
$wsTireRear = $wb.Worksheets.Item("TireRear")
$wsTireRear.Activate()
$wsTireRear.Range("F16").Select()

# --- Brake: re-select B3 (stays the same, but touch it so view state is rewritten) ---
$wsBrake = $wb.Worksheets.Item("Brake")
$wsBrake.Activate()
$wsBrake.Range("B3").Select()

# --- TireFront: becomes the active tab, with new selection B2 ---
# Activated last so it ends up as the workbook's active/visible sheet.
$wsTireFront = $wb.Worksheets.Item("TireFront")
$wsTireFront.Activate()
$wsTireFront.Range("B2").Select()
